# doorstop/core/test/files/exported.xlsx -- "Fixed test failures (still lacking coverage)."
#
# Semantic changes made by the commit:
#   1. E2 ("links" for REQ001): append ":abc123" to the second line, so the
#      cell becomes "SYS001\nSYS002:abc123" instead of "SYS001\nSYS002".
#   2. Column E ("links") is widened from 9.5 to 16.5 (display char width)
#      to accommodate the longer text.
#   3. The whole "reviewed" column (I2:I6) is cleared from boolean FALSE to
#      blank/empty for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "links" value for REQ001 (row 2, column E).
$ws.Range("E2").Value = "SYS001`nSYS002:abc123"

# 2. Widen column E ("links") so the saved OOXML <col> width is 16.5.
#    Excel's COM ColumnWidth property is offset from the stored character
#    width by the standard 5px/6px-per-char padding (~0.8333), so setting
#    ColumnWidth to 15.6666... round-trips to a stored width of 16.5.
$ws.Columns.Item(5).ColumnWidth = 15.666666666666666

# 3. Clear the "reviewed" column (I) for every data row (2-6) to blank.
$ws.Range("I2").Value = ""
$ws.Range("I3").Value = ""
$ws.Range("I4").Value = ""
$ws.Range("I5").Value = ""
$ws.Range("I6").Value = ""
